# Auto-generated update of odds values based on upstream FlashScore diff
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5
$ws.Range("G5").Value = 2.15
$ws.Range("H5").Value = 3.4
$ws.Range("I5").Value = 3.25
$ws.Range("J5").Value = 2.75
$ws.Range("K5").Value = 2.25
$ws.Range("L5").Value = 3.6
$ws.Range("W5").Value = 9.5
$ws.Range("X5").Value = 11
$ws.Range("Y5").Value = 9
$ws.Range("Z5").Value = 21
$ws.Range("AA5").Value = 17
$ws.Range("AB5").Value = 23
$ws.Range("AD5").Value = 6.5
$ws.Range("AE5").Value = 12
$ws.Range("AH5").Value = 12
$ws.Range("AI5").Value = 17
$ws.Range("AJ5").Value = 12
$ws.Range("AK5").Value = 34
$ws.Range("AL5").Value = 23
$ws.Range("AM5").Value = 29
$ws.Range("AN5").Value = 4.33
$ws.Range("AO5").Value = 11
$ws.Range("AQ5").Value = 41
$ws.Range("AS5").Value = 126
$ws.Range("AY5").Value = 17

# Row 6
$ws.Range("G6").Value = 1.57
$ws.Range("H6").Value = 4.33
$ws.Range("I6").Value = 5.25
$ws.Range("J6").Value = 2.1
$ws.Range("Q6").Value = 1.67
$ws.Range("R6").Value = 2.15
$ws.Range("U6").Value = 1.7
$ws.Range("V6").Value = 2.05
$ws.Range("W6").Value = 8.5
$ws.Range("X6").Value = 8.5
$ws.Range("Z6").Value = 12
$ws.Range("AB6").Value = 21
$ws.Range("AE6").Value = 15
$ws.Range("AF6").Value = 41
$ws.Range("AI6").Value = 29
$ws.Range("AK6").Value = 51
$ws.Range("AN6").Value = 3.75
$ws.Range("AQ6").Value = 23
$ws.Range("AX6").Value = 7
$ws.Range("AY6").Value = 26
$ws.Range("AZ6").Value = 29
$ws.Range("BA6").Value = 81
$ws.Range("BC6").Value = 151

# Row 7
$ws.Range("G7").Value = 3.4
$ws.Range("H7").Value = 3.75
$ws.Range("I7").Value = 1.96
$ws.Range("K7").Value = 2.25
$ws.Range("L7").Value = 2.62
$ws.Range("X7").Value = 19
$ws.Range("AD7").Value = 7
$ws.Range("AH7").Value = 9
$ws.Range("AK7").Value = 19
$ws.Range("AN7").Value = 5.5
$ws.Range("AO7").Value = 19
$ws.Range("AP7").Value = 23
$ws.Range("AQ7").Value = 51
$ws.Range("AR7").Value = 67
$ws.Range("AX7").Value = 4.33
$ws.Range("AY7").Value = 11
$ws.Range("BC7").Value = 126

# Row 8
$ws.Range("G8").Value = 6.25
$ws.Range("H8").Value = 4.33
$ws.Range("I8").Value = 1.45
$ws.Range("J8").Value = 5.5
$ws.Range("K8").Value = 2.62
$ws.Range("L8").Value = 1.92
$ws.Range("M8").Value = 1.02
$ws.Range("N8").Value = 19
$ws.Range("O8").Value = 1.13
$ws.Range("P8").Value = 6
$ws.Range("Q8").Value = 1.44
$ws.Range("R8").Value = 2.7
$ws.Range("AC8").Value = 19
$ws.Range("AD8").Value = 9
$ws.Range("AH8").Value = 11
$ws.Range("AY8").Value = 7
$ws.Range("BA8").Value = 19
$ws.Range("BB8").Value = 34

# Row 9
$ws.Range("G9").Value = 1.13
$ws.Range("J9").Value = 1.5

# Row 10
$ws.Range("G10").Value = 1.91
$ws.Range("I10").Value = 3.6
$ws.Range("J10").Value = 2.5
$ws.Range("K10").Value = 2.37
$ws.Range("L10").Value = 3.75
$ws.Range("N10").Value = 17
$ws.Range("Q10").Value = 1.57
$ws.Range("R10").Value = 2.35
$ws.Range("X10").Value = 11
$ws.Range("Y10").Value = 8.5
$ws.Range("Z10").Value = 17
$ws.Range("AH10").Value = 15
$ws.Range("AI10").Value = 21
$ws.Range("AJ10").Value = 13
$ws.Range("AK10").Value = 41
$ws.Range("AL10").Value = 26
$ws.Range("AM10").Value = 29
$ws.Range("AO10").Value = 10
$ws.Range("AY10").Value = 19
$ws.Range("BB10").Value = 67

# Row 11
$ws.Range("G11").Value = 1.95
$ws.Range("I11").Value = 3.4
$ws.Range("J11").Value = 2.4
$ws.Range("L11").Value = 3.5
$ws.Range("N11").Value = 23
$ws.Range("AC11").Value = 23
$ws.Range("AJ11").Value = 13
$ws.Range("AY11").Value = 17
$ws.Range("AZ11").Value = 19
$ws.Range("BA11").Value = 51

# Row 19
$ws.Range("G19").Value = 4.5
$ws.Range("I19").Value = 1.73
$ws.Range("J19").Value = 4.75
$ws.Range("M19").Value = 1.05
$ws.Range("N19").Value = 11
$ws.Range("Q19").Value = 1.85
$ws.Range("R19").Value = 2
$ws.Range("AA19").Value = 34
$ws.Range("AI19").Value = 8.5

# Row 24
$ws.Range("G24").Value = 3.1
$ws.Range("H24").Value = 3.1
$ws.Range("I24").Value = 2.3
$ws.Range("J24").Value = 3.75
$ws.Range("K24").Value = 2
$ws.Range("L24").Value = 3.1
$ws.Range("M24").Value = 1.08
$ws.Range("N24").Value = 8
$ws.Range("O24").Value = 1.4
$ws.Range("P24").Value = 2.75
$ws.Range("Q24").Value = 2.25
$ws.Range("R24").Value = 1.62
$ws.Range("S24").Value = 1.5
$ws.Range("T24").Value = 2.5
$ws.Range("U24").Value = 1.91
$ws.Range("V24").Value = 1.8
$ws.Range("W24").Value = 8.5
$ws.Range("AC24").Value = 8
$ws.Range("AE24").Value = 17
$ws.Range("AL24").Value = 21
$ws.Range("AR24").Value = 101
$ws.Range("AT24").Value = 2.5
$ws.Range("AV24").Value = 67
$ws.Range("AZ24").Value = 26
